# Apply the "X-party 国漫游戏嘉年华03" insertion + follow-up attendee-count
# corrections to the "展览" (sheet1) and "全部类型" (sheet4) worksheets.
#
# Both sheets share the same event list (全部类型 simply has one extra,
# unrelated row earlier on) so the edit is identical on both, just anchored
# at a different insertion row: row 13 on 展览, row 14 on 全部类型.
#
# NOTE: this runtime's PowerShell only reliably binds *positional*
# function arguments, so every helper below takes plain positional params.

function Update-AttendeeCounts($ws, $fixes) {
    foreach ($row in $fixes.Keys) {
        $ws.Cells.Item($row, 6).Value = $fixes[$row]
    }
}

function Insert-NewEvent($ws, $row) {
    # Push the existing row (and everything below it) down by one, without
    # dragging its cell formatting into the new blank row.
    $ws.Rows.Item($row).Insert()
    $ws.Rows.Item($row).ClearFormats()

    # Column A is a plain sequential index (0-based on row 2, so row N -> N-1).
    # Insert() shifts the cell *values* down, but it does not renumber them,
    # so every old index from the insertion point on is now off by one.
    # Walk bottom-up so each cell is only touched once.
    $lastRow = $ws.UsedRange.Rows.Count
    $r = $lastRow
    while ($r -gt $row) {
        $cur = $ws.Cells.Item($r, 1).Value2
        if ($cur -ne $null) {
            $ws.Cells.Item($r, 1).Value = $cur + 1
        }
        $r = $r - 1
    }

    # New row's index cell: give it the same look (bold/centered/bordered)
    # as the column-A cells around it.
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row + 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    # Fill in the new event. B/E look like dates to Excel's text-entry
    # parser, so force them to stay text with a leading quote (same trick
    # Excel itself uses for "numbers stored as text").
    $ws.Cells.Item($row, 2).Value = "'2024.04.05"
    $ws.Cells.Item($row, 3).Value = "苏州·X-party 国漫游戏嘉年华03"
    $ws.Cells.Item($row, 4).Value = "秋枫街与开平路交叉口西南角 爱琴海购物中心"
    $ws.Cells.Item($row, 5).Value = "'2024.04.05 10:00-04.06 17:00"
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 48
    $ws.Cells.Item($row, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82042"
    $ws.Cells.Item($row, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/WaQk4nUt1708679999084.jpeg"
}

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$fixes1 = @{ 3 = 87; 4 = 107; 5 = 43; 7 = 2707; 9 = 957; 10 = 275; 12 = 10321 }
Update-AttendeeCounts $ws1 $fixes1
Insert-NewEvent $ws1 13

# ---- Sheet "全部类型" -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$fixes4 = @{ 3 = 87; 4 = 107; 5 = 43; 7 = 2707; 10 = 957; 11 = 275; 13 = 10321 }
Update-AttendeeCounts $ws4 $fixes4
Insert-NewEvent $ws4 14
